$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix data point: Palo Alto 2019 budget (C6) was 95000000, should be 9500000
$ws.Range("C6").Value = 9500000

# Update selection to reflect active cell C6
$ws.Range("C6").Select()
